# feat: Add MAL agg packages for v 2.35, 2.36; update for v 2.33, 2.34
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Package info": refresh version/build metadata and append a Name row
# ---------------------------------------------------------------------------
$wsInfo = $wb.Worksheets.Item("Package info")

$wsInfo.Range("B4").Value = "1.1.0"
$wsInfo.Range("B5").Value = "2.33.9"

$wsInfo.Range("A6").Value = "DHIS2 build"
$wsInfo.Range("B6").Value = "58094d2"

$wsInfo.Range("A7").Value = "Last updated"
$wsInfo.Range("B7").Value = "20210520T090044"

$wsInfo.Range("A8").Value = "Name"
$wsInfo.Range("B8").Value = "MAL-MMG_CUSTOM_V1.1.0_2.33.9-en"
# give the new row the same banding style as the other "odd" data rows (row 6)
$wsInfo.Range("A6:B6").Copy()
$wsInfo.Range("A8:B8").PasteSpecial(-4122)

# narrow column B to match the new, shorter values
$wsInfo.Columns.Item(2).ColumnWidth = 32.83

# ---------------------------------------------------------------------------
# Sheet "dataElements": the six data-element rows were reordered
# ---------------------------------------------------------------------------
$wsDE = $wb.Worksheets.Item("dataElements")

$deRows = @(
    @("MAL - Migrant and mobile population (MMP) positive", "MMP positive", "MAL_MM_POP_POS", "Number of migrant and mobile population (MMP) that are positive with either microscopy and RDT", "bjDvmb4bfuf", "2019-10-20", "CWHBMa4nC9J"),
    @("MAL - Malaria tested from cross-borders", "Tested from cross-borders", "MAL_TEST_CROSS_BORDERS", "Number of suspected cases tested with either microscopy or RDT in administrative area bordering international border", "VkQPxB6VdoG", "2019-10-20", "CxI1FHE4oEh"),
    @("MAL - Migrant and mobile population (MMP) followed up for 14 days", "MMP followed up for 14 days", "MAL_MM_POP_FOLLO_UP_FOR_14D", "Number of migrant and mobile population (MMP) that are are followed-up for 14 days (with testing using microscopy or RDT at intervals)", "bjDvmb4bfuf", "2019-10-20", "kdMT3AuDzj1"),
    @("MAL - Migrant and mobile population (MMP) tested", "MMP tested", "MAL_MM_POP_TEST", "Number of migrant and mobile population (MMP) suspected and tested with either microscopy and RDT", "bjDvmb4bfuf", "2019-10-20", "S3AqkeU4DET"),
    @("MAL - Malaria positive from cross-borders", "Positive from cross-borders", "MAL_POS_CROSS_BORDERS", "Number of positive cases with either microscopy or RDT in administrative area bordering international border", "VkQPxB6VdoG", "2019-10-20", "UwaQ0MJzXBz"),
    @("MAL - Malaria positive from cross-borders followed for 14 days", "Positive from cross-borders followed for 14 days", "MAL_POS_CROSS_BORDERS_FOLLO_14D", "Number of suspected malaria cases positive with either microscopy or RDT in administrative area bordering international border", "bjDvmb4bfuf", "2019-10-20", "wAHUeGPbH9A")
)

# Column F ("Last updated") is "2019-10-20" for every single row both before
# and after the edit, so it is intentionally left untouched here.
$r = 2
foreach ($row in $deRows) {
    $wsDE.Range("A$r").Value = $row[0]
    $wsDE.Range("B$r").Value = $row[1]
    $wsDE.Range("C$r").Value = $row[2]
    $wsDE.Range("D$r").Value = $row[3]
    $wsDE.Range("E$r").Value = $row[4]
    $wsDE.Range("G$r").Value = $row[6]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Sheet "dataElementGroups": same six data elements, new order in column B
# ---------------------------------------------------------------------------
$wsDEG = $wb.Worksheets.Item("dataElementGroups")

$degRows = @(
    "MAL - Migrant and mobile population (MMP) positive",
    "MAL - Malaria tested from cross-borders",
    "MAL - Migrant and mobile population (MMP) followed up for 14 days",
    "MAL - Migrant and mobile population (MMP) tested",
    "MAL - Malaria positive from cross-borders",
    "MAL - Malaria positive from cross-borders followed for 14 days"
)

$r = 2
foreach ($name in $degRows) {
    $wsDEG.Range("B$r").Value = $name
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Sheet "userGroups": "Malaria access" and "Malaria data capture" swap places,
# and every row's "Last updated" date is refreshed
# ---------------------------------------------------------------------------
$wsUG = $wb.Worksheets.Item("userGroups")

# Force the date-shaped values to stay plain text (matches the source file's
# t="str" cells) instead of being auto-parsed into date serials.
$wsUG.Range("B2:B4").NumberFormat = "@"

$wsUG.Range("A2").Value = "Malaria data capture"
$wsUG.Range("B2").Value = "2021-05-20"
$wsUG.Range("C2").Value = "fRSrUJ6SMGH"

$wsUG.Range("A3").Value = "Malaria admin"
$wsUG.Range("B3").Value = "2021-05-20"
$wsUG.Range("C3").Value = "suMb19wGXPR"

$wsUG.Range("A4").Value = "Malaria access"
$wsUG.Range("B4").Value = "2021-05-20"
$wsUG.Range("C4").Value = "ZXEVDM9XRea"
